# issue237/fill-list.xlsx demo-data refresh
#
# - Sheet1!A1 shared string "字符串标题" -> "StringTitle"
# - Sheet1!B1 shared string "数字标题"   -> "DateTitle"
# - Sheet1 active selection moves from A8 to B6
# - (bookViews windowHeight is a cosmetic last-saved-UI-state value that
#   isn't exposed on the Excel object model; left untouched.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the demo header labels (find & replace keeps cell formatting/style
# intact while swapping the underlying text).
$ws.Cells.Replace("字符串标题", "StringTitle")
$ws.Cells.Replace("数字标题", "DateTitle")

# Move the sheet's active selection to B6.
$ws.Range("B6").Select()
